# Fix LBNRIND label (was mistakenly "LBNDIND") in the header row of the
# reference-ranges lookup table, and clear the now-unused "LBNDIND" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the typo'd column header: "LBNDIND" -> "LBNRIND"
$ws.Range("F1").Value = "LBNRIND"

# Replicate the saved cursor/selection position from the author's session.
$ws.Range("I9").Select()
